# Update the "Lists" sheet's Browser list:
#  - row 8 (previously "Headless") becomes "Chrome-Headless"
#  - a new row 9 "FireFox-Headless" is appended
# Also move the sheet's selection to F3 to match the saved view state,
# without disturbing which sheet/cell is active on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws = $wb.Worksheets.Item("Lists")

$ws.Range("A8").Value = "Chrome-Headless"
$ws.Range("A9").Value = "FireFox-Headless"

$ws.Range("F3").Select()

# Restore Sheet1 as the active/selected sheet (selecting a range on
# "Lists" above implicitly activates that sheet).
$ws1.Activate()
